$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the "Normal" style to the existing data range so a new cellXfs
# entry (applyFont="true", same font) is created and used by A1:B6.
$ws.Range("A1:B6").Style = "Normal"

# Add new row of data: A7 = 6 (no B7 value)
$ws.Range("A7").Value = 6

# Move the selection like the author did (A8)
$ws.Range("A8").Select()
